# Appends the "fiche de personnage" stat block (skills table, Attributs,
# Inventaires, Sortileges sections) as a run of new paragraphs after the
# two existing summary tables, right before the section properties.
$d = $word.ActiveDocument

# The document's final (sentinel) paragraph sits immediately after the
# second table, right before the body's sectPr - that's where the new
# content needs to land.
$p = $d.Paragraphs.Last
$r = $p.Range

# Build the new paragraphs as a literal WordprocessingML fragment (tabs as
# real <w:tab/> run children, not literal tab characters) and splice it in
# with InsertXML, which replaces the target range's contents in place.
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr/><w:r><w:tab/><w:t xml:space="preserve">Adresse : 5</w:t></w:r><w:r><w:tab/><w:tab/><w:tab/><w:tab/><w:t xml:space="preserve">Explosifs : 0</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Agilité : 0</w:t></w:r><w:r><w:tab/><w:tab/><w:tab/><w:tab/><w:t xml:space="preserve">Force : 0</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Animale : 0</w:t></w:r><w:r><w:tab/><w:tab/><w:tab/><w:tab/><w:t xml:space="preserve">Intimidation : 0</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Artisanat : 0</w:t></w:r><w:r><w:tab/><w:tab/><w:tab/><w:tab/><w:t xml:space="preserve">Langages : 0</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Botanique : 0</w:t></w:r><w:r><w:tab/><w:tab/><w:tab/><w:tab/><w:t xml:space="preserve">Mécanique : 5</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Connaissances géographiques : 0</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve">Médecine : 0</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Connaissances historiques : 0</w:t></w:r><w:r><w:tab/><w:tab/><w:t xml:space="preserve">Natation : 0</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Connaissances magiques : 0</w:t></w:r><w:r><w:tab/><w:tab/><w:t xml:space="preserve">Perception : 0</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Connaissances religieuse : 0</w:t></w:r><w:r><w:tab/><w:tab/><w:t xml:space="preserve">Perspicacité : 0</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Crochetage : 0</w:t></w:r><w:r><w:tab/><w:tab/><w:tab/><w:tab/><w:t xml:space="preserve">Persuasion : 0</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Diplomatie : 0</w:t></w:r><w:r><w:tab/><w:tab/><w:tab/><w:tab/><w:t xml:space="preserve">Psyché : 0</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Discrétion : 0</w:t></w:r><w:r><w:tab/><w:tab/><w:tab/><w:tab/><w:t xml:space="preserve">Réflexes : 0</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Endurance : 0</w:t></w:r><w:r><w:tab/><w:tab/><w:tab/><w:tab/><w:t xml:space="preserve">Vigueur : 0</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Escalade : 0</w:t></w:r><w:r><w:tab/><w:tab/><w:tab/><w:tab/><w:t xml:space="preserve">Volonté : 0</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr/></w:pPr></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr/><w:r><w:t xml:space="preserve">Attributs : </w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">Avantage du terrain: sur x terrain(s), la créature n'a pas de malus, Fin limier: plafond supplémentaire de 5% dans une des compétences techniques</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr/><w:r><w:t xml:space="preserve">Inventaires : </w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr/></w:pPr><w:r><w:t xml:space="preserve">Sortilèges : </w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr/></w:pPr></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr/></w:p>
'@

$r.InsertXML($xml)
